$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force "Text" number format on cells whose new price value would otherwise
# be auto-parsed by Excel as a number (losing the exact textual representation,
# e.g. trailing zeros like "1.00" or "5.70").
$textCells = @(
    "D5",
    "D6",
    "D8",
    "D10",
    "D13",
    "D18",
    "D19",
    "D20",
    "D21",
    "D23",
    "D25",
    "D27",
    "D29",
    "D30",
    "D31",
    "D32",
    "D36",
    "D39",
    "D40",
    "D42",
    "D43",
    "D44",
    "D45",
    "D46",
    "D48",
    "D49",
    "D50",
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated Price (D) / Volume(1h) (E) values scraped by the Action run.
$ws.Range("D2").Value = "59.700.11"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "2.403.10"
$ws.Range("E3").Value = "  -0.34%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "549.90"
$ws.Range("E5").Value = "  -0.48%  "
$ws.Range("D6").Value = "136.65"
$ws.Range("E6").Value = "  -0.32%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "0.592"
$ws.Range("E8").Value = "  +4.23%  "
$ws.Range("E9").Value = "  -2.32%  "
$ws.Range("D10").Value = "5.70"
$ws.Range("E10").Value = "  -1.69%  "
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("E12").Value = "  -1.91%  "
$ws.Range("D13").Value = "25.22"
$ws.Range("E13").Value = "  +2.53%  "
$ws.Range("D14").Value = "2.831.33"
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("D15").Value = "59.656.11"
$ws.Range("E15").Value = "  +0.35%  "
$ws.Range("E16").Value = "  -1.77%  "
$ws.Range("D17").Value = "2.390.38"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("D18").Value = "11.29"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").Value = "4.40"
$ws.Range("E19").Value = "  -1.07%  "
$ws.Range("D20").Value = "328.04"
$ws.Range("E20").Value = "  -1.96%  "
$ws.Range("D21").Value = "6.65"
$ws.Range("E21").Value = "  -3.98%  "
$ws.Range("E22").Value = "  -0.03%  "
$ws.Range("D23").Value = "66.43"
$ws.Range("E23").Value = "  +3.25%  "
$ws.Range("E24").Value = "  +1.54%  "
$ws.Range("D25").Value = "8.62"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").Value = "1.37"
$ws.Range("E27").Value = "  -1.23%  "
$ws.Range("D28").Value = "0.0₃0768"
$ws.Range("E28").Value = "  -2.30%  "
$ws.Range("D29").Value = "1.77"
$ws.Range("E29").Value = "  -1.97%  "
$ws.Range("D30").Value = "167.99"
$ws.Range("E30").Value = "  -1.43%  "
$ws.Range("D31").Value = "6.06"
$ws.Range("E31").Value = "  -3.44%  "
$ws.Range("D32").Value = "18.58"
$ws.Range("E32").Value = "  -0.52%  "
$ws.Range("E33").Value = "  -2.12%  "
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.09%  "
$ws.Range("E37").Value = "  -2.18%  "
$ws.Range("E38").Value = "  -1.87%  "
$ws.Range("D39").Value = "312.93"
$ws.Range("E39").Value = "  +2.10%  "
$ws.Range("D40").Value = "0.407"
$ws.Range("E40").Value = "  -3.17%  "
$ws.Range("E41").Value = "  -2.31%  "
$ws.Range("D42").Value = "138.71"
$ws.Range("E42").Value = "  -2.81%  "
$ws.Range("D43").Value = "0.0969"
$ws.Range("E43").Value = "  +0.43%  "
$ws.Range("D44").Value = "0.0515"
$ws.Range("E44").Value = "  -1.57%  "
$ws.Range("D45").Value = "19.46"
$ws.Range("E45").Value = "  +1.45%  "
$ws.Range("D46").Value = "0.577"
$ws.Range("E46").Value = "  +1.00%  "
$ws.Range("E47").Value = "  -1.62%  "
$ws.Range("D48").Value = "0.386"
$ws.Range("E48").Value = "  -5.80%  "
$ws.Range("D49").Value = "17.54"
$ws.Range("E49").Value = "  -2.14%  "
$ws.Range("D50").Value = "11.05"
$ws.Range("E50").Value = "  +0.08%  "
$ws.Range("E51").Value = "  -3.07%  "
